$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A2').Value2 = 111896654
$ws.Range('B2').Value2 = 89183
$ws.Range('E2').Value2 = 3215
$ws.Range('F2').Value2 = 'Rödgul trumpetsvamp'
$ws.Range('G2').Value2 = 'Craterellus lutescens'
$ws.Range('H2').Value2 = '(Fr.) Fr.'
$ws.Range('Q2').Value2 = 575072.6962527435
$ws.Range('R2').Value2 = 6703421.833381963
$ws.Range('A3').Value2 = 111896636
$ws.Range('Q3').Value2 = 575108.85141061
$ws.Range('R3').Value2 = 6703418.142308297
$ws.Range('A4').Value2 = 111896642
$ws.Range('Q4').Value2 = 575014.1091647458
$ws.Range('R4').Value2 = 6703387.066676207
$ws.Range('A6').Value2 = 111896640
$ws.Range('B6').Value2 = 90332
$ws.Range('E6').Value2 = 4769
$ws.Range('F6').Value2 = 'Svavelriska'
$ws.Range('G6').Value2 = 'Lactarius scrobiculatus'
$ws.Range('H6').Value2 = '(Scop.:Fr.) Fr.'
$ws.Range('Q6').Value2 = 575025.3556637274
$ws.Range('R6').Value2 = 6703369.042946251
$ws.Range('A7').Value2 = 111896637
$ws.Range('Q7').Value2 = 575088.0587098968
$ws.Range('R7').Value2 = 6703396.00058554
$ws.Range('A8').Value2 = 111884471
$ws.Range('B8').Value2 = 88899
$ws.Range('D8').Value2 = 'NT'
$ws.Range('E8').Value2 = 3286
$ws.Range('F8').Value2 = 'Flattoppad klubbsvamp'
$ws.Range('G8').Value2 = 'Clavariadelphus truncatus'
$ws.Range('H8').Value2 = '(Quél.) Donk'
$ws.Range('Q8').Value2 = 575020.8210917887
$ws.Range('R8').Value2 = 6703397.074168184
$ws.Range('A9').Value2 = 111883983
$ws.Range('K9').Value2 = "'"
$ws.Range('K9').Style = "Normal"
$ws.Range('P9').Value2 = 'Kalkberget (Kalkberget), Gstr'
$ws.Range('Q9').Value2 = 575058.3527020445
$ws.Range('R9').Value2 = 6703446.206921679
$ws.Range('AW9').Value2 = 'Patric Engfeldt'
$ws.Range('AX9').Value2 = 'Patric Engfeldt'
$ws.Range('A10').Value2 = 111896653
$ws.Range('Q10').Value2 = 575075.050630242
$ws.Range('R10').Value2 = 6703403.625642136
$ws.Range('A11').Value2 = 111896633
$ws.Range('Q11').Value2 = 575100.4050603262
$ws.Range('R11').Value2 = 6703444.118284944
$ws.Range('A12').Value2 = 111896635
$ws.Range('B12').Value2 = 90332
$ws.Range('E12').Value2 = 4769
$ws.Range('F12').Value2 = 'Svavelriska'
$ws.Range('G12').Value2 = 'Lactarius scrobiculatus'
$ws.Range('H12').Value2 = '(Scop.:Fr.) Fr.'
$ws.Range('Q12').Value2 = 575037.2974304935
$ws.Range('R12').Value2 = 6703389.027347369
$ws.Range('AF12').ClearContents()
$ws.Range('A13').Value2 = 111896643
$ws.Range('Q13').Value2 = 575038.7114136803
$ws.Range('R13').Value2 = 6703416.194821274
$ws.Range('A14').Value2 = 111896690
$ws.Range('B14').Value2 = 90687
$ws.Range('E14').Value2 = 5964
$ws.Range('F14').Value2 = 'Fjällig taggsvamp s.str.'
$ws.Range('G14').Value2 = 'Sarcodon imbricatus s.str.'
$ws.Range('H14').Value2 = '(L.:Fr.) P.Karst.'
$ws.Range('Q14').Value2 = 575060.2881161601
$ws.Range('R14').Value2 = 6703376.67477417
$ws.Range('AF14').Value2 = "'"
$ws.Range('AF14').Style = "Normal"
$ws.Range('A15').Value2 = 111884093
$ws.Range('B15').Value2 = 98535
$ws.Range('D15').Value2 = 'LC'
$ws.Range('E15').Value2 = 222498
$ws.Range('F15').Value2 = 'Blåsippa'
$ws.Range('G15').Value2 = 'Hepatica nobilis'
$ws.Range('H15').Value2 = 'Schreb.'
$ws.Range('P15').Value2 = 'Kopparåsen (Kopparåsen), Gstr'
$ws.Range('Q15').Value2 = 575065.9914513066
$ws.Range('R15').Value2 = 6703387.648325931
$ws.Range('A16').Value2 = 111884133
$ws.Range('B16').Value2 = 88899
$ws.Range('D16').Value2 = 'NT'
$ws.Range('E16').Value2 = 3286
$ws.Range('F16').Value2 = 'Flattoppad klubbsvamp'
$ws.Range('G16').Value2 = 'Clavariadelphus truncatus'
$ws.Range('H16').Value2 = '(Quél.) Donk'
$ws.Range('K16').Value2 = "'"
$ws.Range('K16').Style = "Normal"
$ws.Range('P16').Value2 = 'Kalkberget (Kalkberget), Gstr'
$ws.Range('Q16').Value2 = 575059.034285416
$ws.Range('R16').Value2 = 6703389.477814267
$ws.Range('AW16').Value2 = 'Patric Engfeldt'
$ws.Range('AX16').Value2 = 'Patric Engfeldt'
$ws.Range('A17').Value2 = 111896641
$ws.Range('Q17').Value2 = 575021.3626164712
$ws.Range('R17').Value2 = 6703370.933926445
$ws.Range('A18').Value2 = 111896652
$ws.Range('B18').Value2 = 89183
$ws.Range('D18').Value2 = 'LC'
$ws.Range('E18').Value2 = 3215
$ws.Range('F18').Value2 = 'Rödgul trumpetsvamp'
$ws.Range('G18').Value2 = 'Craterellus lutescens'
$ws.Range('H18').Value2 = '(Fr.) Fr.'
$ws.Range('K18').ClearContents()
$ws.Range('P18').Value2 = 'Kratte masugn, Gstr'
$ws.Range('Q18').Value2 = 575066.556649723
$ws.Range('R18').Value2 = 6703455.751857814
$ws.Range('AW18').Value2 = 'Philipp Weiss'
$ws.Range('AX18').Value2 = 'Philipp Weiss'
$ws.Range('A19').Value2 = 111896655
$ws.Range('B19').Value2 = 89183
$ws.Range('E19').Value2 = 3215
$ws.Range('F19').Value2 = 'Rödgul trumpetsvamp'
$ws.Range('G19').Value2 = 'Craterellus lutescens'
$ws.Range('H19').Value2 = '(Fr.) Fr.'
$ws.Range('K19').ClearContents()
$ws.Range('P19').Value2 = 'Kratte masugn, Gstr'
$ws.Range('Q19').Value2 = 575104.6742508161
$ws.Range('R19').Value2 = 6703428.910891063
$ws.Range('AW19').Value2 = 'Philipp Weiss'
$ws.Range('AX19').Value2 = 'Philipp Weiss'
$ws.Range('A20').Value2 = 111896638
$ws.Range('Q20').Value2 = 575087.1320314853
$ws.Range('R20').Value2 = 6703393.020834555
$ws.Range('A21').Value2 = 111896644
$ws.Range('Q21').Value2 = 575036.4083237475
$ws.Range('R21').Value2 = 6703431.936489306
$ws.Range('A22').Value2 = 111896639
$ws.Range('B22').Value2 = 90332
$ws.Range('E22').Value2 = 4769
$ws.Range('F22').Value2 = 'Svavelriska'
$ws.Range('G22').Value2 = 'Lactarius scrobiculatus'
$ws.Range('H22').Value2 = '(Scop.:Fr.) Fr.'
$ws.Range('Q22').Value2 = 575089.384229039
$ws.Range('R22').Value2 = 6703379.745088123
